$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 156.5
$ws.Range("I11").Value = 156.5
$ws.Range("K11").Value = 156.5
$ws.Range("M11").Value = -16.5
$ws.Range("H33").Value = 193.26923
$ws.Range("I33").Value = 197.77272
$ws.Range("J33").Value = 168.5
$ws.Range("K33").Value = 197.77272
$ws.Range("L33").Value = 168.5
$ws.Range("M33").Value = 31.22728000000001
$ws.Range("N33").Value = -626.5
$ws.Range("H44").Value = 50050
$ws.Range("J44").Value = 50050
$ws.Range("L44").Value = 50050
$ws.Range("N44").Value = -50974
$ws.Range("H86").Value = 4389464
$ws.Range("I86").Value = 2498.1667
$ws.Range("J86").Value = 8776429
$ws.Range("K86").Value = 2498.1667
$ws.Range("L86").Value = 8776429
$ws.Range("M86").Value = -1375.1667
$ws.Range("N86").Value = -8778675
$ws.Range("H89").Value = 4389464
$ws.Range("I89").Value = 2498.1667
$ws.Range("J89").Value = 8776429
$ws.Range("K89").Value = 12490.8335
$ws.Range("L89").Value = 43882145
$ws.Range("M89").Value = -6874.833500000001
$ws.Range("N89").Value = -43893377
$ws.Range("H96").Value = 884.2857
$ws.Range("J96").Value = 902.25
$ws.Range("L96").Value = 2706.75
$ws.Range("N96").Value = -5452.75
$ws.Range("H101").Value = 444.16666
$ws.Range("I101").Value = 353.75
$ws.Range("K101").Value = 1061.25
$ws.Range("M101").Value = 560.75
$ws.Range("H111").Value = 3857.923
$ws.Range("J111").Value = 4499.5
$ws.Range("L111").Value = 13498.5
$ws.Range("N111").Value = -19632.5
$ws.Range("H118").Value = 13935.875
$ws.Range("I118").Value = 18396.334
$ws.Range("K118").Value = 55189.00199999999
$ws.Range("M118").Value = -53532.00199999999
$ws.Range("H126").Value = 28333.2
$ws.Range("J126").Value = 28333.2
$ws.Range("L126").Value = 28333.2
$ws.Range("N126").Value = -38213.2
$ws.Range("H127").Value = 9140.286
$ws.Range("I127").Value = 2401.2856
$ws.Range("J127").Value = 15879.286
$ws.Range("K127").Value = 7203.8568
$ws.Range("L127").Value = 47637.858
$ws.Range("M127").Value = -2243.8568
$ws.Range("N127").Value = -57557.858
$ws.Range("H132").Value = 20766.643
$ws.Range("I132").Value = 1994.174
$ws.Range("K132").Value = 5982.522
$ws.Range("M132").Value = -3452.522
$ws.Range("H135").Value = 13159956
$ws.Range("I135").Value = 22729532
$ws.Range("K135").Value = 204565788
$ws.Range("M135").Value = -204563253
$ws.Range("H137").Value = 5723520.5
$ws.Range("I137").Value = 10013769
$ws.Range("J137").Value = 3189.4
$ws.Range("K137").Value = 30041307
$ws.Range("L137").Value = 9568.200000000001
$ws.Range("M137").Value = -30038757
$ws.Range("N137").Value = -14668.2
$ws.Range("H138").Value = 4437.475
$ws.Range("I138").Value = 3304.2856
$ws.Range("J138").Value = 5047.654
$ws.Range("K138").Value = 9912.856800000001
$ws.Range("L138").Value = 15142.962
$ws.Range("M138").Value = -4772.856800000001
$ws.Range("N138").Value = -25422.962

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2376.5417
$ws.Range("I61").Value = 1856.5238
$ws.Range("K61").Value = 1856.5238
$ws.Range("M61").Value = -1644.5238
$ws.Range("H74").Value = 1537.8667
$ws.Range("I74").Value = 1201.25
$ws.Range("K74").Value = 1201.25
$ws.Range("M74").Value = -327.25
$ws.Range("H77").Value = 1537.8667
$ws.Range("I77").Value = 1201.25
$ws.Range("K77").Value = 6006.25
$ws.Range("M77").Value = -1638.25
$ws.Range("H122").Value = 4918.514
$ws.Range("I122").Value = 4450.7827
$ws.Range("K122").Value = 13352.3481
$ws.Range("M122").Value = -10902.3481
$ws.Range("H136").Value = 2376.5417
$ws.Range("I136").Value = 1856.5238
$ws.Range("K136").Value = 5569.5714
$ws.Range("M136").Value = -3019.5714

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1717.5834
$ws.Range("I105").Value = 1328.2727
$ws.Range("K105").Value = 1328.2727
$ws.Range("M105").Value = 418.7273
$ws.Range("H137").Value = 59571.43
$ws.Range("J137").Value = 59571.43
$ws.Range("L137").Value = 59571.43
$ws.Range("N137").Value = -69771.42999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5574.5405
$ws.Range("I31").Value = 6059.875
$ws.Range("J31").Value = 4678.5386
$ws.Range("K31").Value = 6059.875
$ws.Range("L31").Value = 4678.5386
$ws.Range("M31").Value = -5764.875
$ws.Range("N31").Value = -5268.5386
$ws.Range("H34").Value = 5574.5405
$ws.Range("I34").Value = 6059.875
$ws.Range("J34").Value = 4678.5386
$ws.Range("K34").Value = 6059.875
$ws.Range("L34").Value = 4678.5386
$ws.Range("M34").Value = -5857.875
$ws.Range("N34").Value = -5082.5386
$ws.Range("H58").Value = 2052.1765
$ws.Range("I58").Value = 1294.8462
$ws.Range("K58").Value = 1294.8462
$ws.Range("M58").Value = -1091.8462
$ws.Range("H105").Value = 2782.3333
$ws.Range("I105").Value = 2599.4
$ws.Range("K105").Value = 2599.4
$ws.Range("M105").Value = -852.4000000000001
$ws.Range("H125").Value = 72081.5
$ws.Range("J125").Value = 72081.5
$ws.Range("L125").Value = 72081.5
$ws.Range("N125").Value = -77001.5
$ws.Range("H136").Value = 2052.1765
$ws.Range("I136").Value = 1294.8462
$ws.Range("K136").Value = 3884.5386
$ws.Range("M136").Value = -1334.5386

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 999500
$ws.Range("I128").Value = 999500
$ws.Range("K128").Value = 2998500
$ws.Range("M128").Value = -2993520
$ws.Range("H137").Value = 1302.6
$ws.Range("J137").Value = 2033
$ws.Range("L137").Value = 6099
$ws.Range("N137").Value = -16299

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4551351
$ws.Range("I97").Value = 743.13336
$ws.Range("K97").Value = 743.13336
$ws.Range("M97").Value = -247.13336
$ws.Range("H102").Value = 3579.5264
$ws.Range("I102").Value = 2963.375
$ws.Range("K102").Value = 2963.375
$ws.Range("M102").Value = -1341.375
$ws.Range("H126").Value = 2583.3333
$ws.Range("I126").Value = 2583.3333
$ws.Range("K126").Value = 7749.999899999999
$ws.Range("M126").Value = -5279.999899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 750
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -405
$ws.Range("N22").Value = -1490
$ws.Range("H27").Value = 750
$ws.Range("I27").Value = 700
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 700
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -593
$ws.Range("N27").Value = -1114
$ws.Range("H40").Value = 3997.25
$ws.Range("I40").Value = 4179.4546
$ws.Range("J40").Value = 1993
$ws.Range("K40").Value = 4179.4546
$ws.Range("L40").Value = 1993
$ws.Range("M40").Value = -4043.4546
$ws.Range("N40").Value = -2265
$ws.Range("H46").Value = 1813.84
$ws.Range("I46").Value = 591.25
$ws.Range("J46").Value = 2046.7142
$ws.Range("K46").Value = 591.25
$ws.Range("L46").Value = 2046.7142
$ws.Range("M46").Value = -403.25
$ws.Range("N46").Value = -2422.7142
$ws.Range("H61").Value = 3396.8
$ws.Range("I61").Value = 3372.25
$ws.Range("J61").Value = 3495
$ws.Range("K61").Value = 3372.25
$ws.Range("L61").Value = 3495
$ws.Range("M61").Value = -3170.25
$ws.Range("N61").Value = -3899
$ws.Range("H100").Value = 83687.69500000001
$ws.Range("H105").Value = 300000
$ws.Range("J105").Value = 300000
$ws.Range("L105").Value = 300000
$ws.Range("N105").Value = -306988
$ws.Range("H113").Value = 3396.8
$ws.Range("I113").Value = 3372.25
$ws.Range("J113").Value = 3495
$ws.Range("K113").Value = 3372.25
$ws.Range("L113").Value = 3495
$ws.Range("M113").Value = -1202.25
$ws.Range("N113").Value = -7835
$ws.Range("H122").Value = 4600.6
$ws.Range("J122").Value = 6499.5
$ws.Range("L122").Value = 19498.5
$ws.Range("N122").Value = -24398.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 23494614
$ws.Range("I81").Value = 14286725
$ws.Range("K81").Value = 28573450
$ws.Range("M81").Value = -28572389
$ws.Range("H84").Value = 23494614
$ws.Range("I84").Value = 14286725
$ws.Range("K84").Value = 142867250
$ws.Range("M84").Value = -142861946
$ws.Range("H96").Value = 102602
$ws.Range("I96").Value = 250500
$ws.Range("J96").Value = 4003.3333
$ws.Range("K96").Value = 250500
$ws.Range("L96").Value = 4003.3333
$ws.Range("M96").Value = -249127
$ws.Range("N96").Value = -6749.3333
$ws.Range("H113").Value = 1365.1666
$ws.Range("I113").Value = 1365.1666
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4095.4998
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1925.4998
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 3879.8
$ws.Range("I126").Value = 3624.75
$ws.Range("K126").Value = 10874.25
$ws.Range("M126").Value = -8404.25
$ws.Range("H136").Value = 3820.718
$ws.Range("I136").Value = 2953.6333
$ws.Range("K136").Value = 8860.8999
$ws.Range("M136").Value = -6310.8999
